# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig -- update FHIR
# StructureDefinition metadata (Version/Date/Publisher/Jurisdiction) and
# the Extension's Short/Definition text on the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# The old sheet had a duplicated "Contact" / "No display for ContactDetail"
# row (rows 10 & 11). Remove the duplicate row so everything below shifts
# up by one (this also tightens the sheet dimension from B21 to B20).
$meta.Rows.Item(11).Delete()

# Version bump.
$meta.Range("B3").Value = "6.0.0"

# New IG build date.
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value.
$meta.Range("B9").Value = "Alvearie Team"

# The old duplicate "Contact" row (now row 10, after the delete above)
# becomes a new "Jurisdiction" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root extension row: Short / Definition text updated to describe this
# specific extension instead of the generic boilerplate.
$elements.Range("K2").Value = "Employee Business Unit"
$elements.Range("L2").Value = "Business unit of the subscriber as reported on the claim"
